$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting the existing rows 93:139 down to 94:140.
$ws.Rows(93).Insert()

# Populate the newly inserted row 93 with this week's record (same market/region/
# product metadata, new date + price figures).
$ws.Cells.Item(93, 1).Value = 11
$ws.Cells.Item(93, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(93, 3).Value = "Bíobío"
$ws.Cells.Item(93, 4).Value = 45202
$ws.Cells.Item(93, 5).Value = 8
$ws.Cells.Item(93, 6).Value = 100112037
$ws.Cells.Item(93, 7).Value = "Cebollín"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 50
$ws.Cells.Item(93, 11).Value = 4000
$ws.Cells.Item(93, 12).Value = 4500
$ws.Cells.Item(93, 13).Value = 4300
$ws.Cells.Item(93, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(93, 15).Value = "Región Metropolitana"
$ws.Cells.Item(93, 16).Value = 119
$ws.Cells.Item(93, 17).Value = 36
$ws.Cells.Item(93, 18).Value = "Hortaliza"
